# Apply the "cryptos list" refresh (Tue Sep 12 12:28:17 UTC 2023, GitHub Actions).
# Most rows keep the same coin in place and only refresh the Price (D) and/or
# Volume(1h) (E) columns. Rows 37-40 additionally reshuffle which coin sits on
# which row (VeChain/PaxDollar/ImmutableX/MXToken rotate), so those four rows
# get every column (B/C/D/E) rewritten explicitly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet stores every Price cell as literal text (e.g. "213.04"), even
# though it looks like a plain number. Excel's COM Value setter auto-detects
# "clean" decimal text and coerces it to a real number (dropping formatting like
# trailing zeros, e.g. "6.50" -> 6.5). Pre-marking those specific Price cells as
# Text keeps the assignment a literal string, matching the source data shape.
# (Multi-dot values like "26.211.91" are never auto-numeric, so they are left alone.)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7:D9").NumberFormat = "@"
$ws.Range("D14:D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D20:D25").NumberFormat = "@"
$ws.Range("D28:D30").NumberFormat = "@"
$ws.Range("D32:D34").NumberFormat = "@"
$ws.Range("D37:D42").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"

# --- Row 2: Bitcoin ---
$ws.Range("D2").Value = "26.211.91"
$ws.Range("E2").Value = "  +1.62%  "

# --- Row 3: Ethereum ---
$ws.Range("D3").Value = "1.612.02"
$ws.Range("E3").Value = "  +0.94%  "

# --- Row 4: TetherUSD ---
$ws.Range("E4").Value = "  -0.46%  "

# --- Row 5: BNB ---
$ws.Range("D5").Value = "212.77"
$ws.Range("E5").Value = "  +1.87%  "

# --- Row 7: XRP ---
$ws.Range("D7").Value = "0.483"
$ws.Range("E7").Value = "  +1.30%  "

# --- Row 8: Cardano ---
$ws.Range("D8").Value = "0.249"
$ws.Range("E8").Value = "  +1.44%  "

# --- Row 9: Dogecoin ---
$ws.Range("D9").Value = "0.0620"
$ws.Range("E9").Value = "  +1.87%  "

# --- Row 10: Solana ---
$ws.Range("E10").Value = "  +3.20%  "

# --- Row 11: TRON ---
$ws.Range("E11").Value = "  +1.77%  "

# --- Row 12: WrappedliquidstakedEther2.0 ---
$ws.Range("D12").Value = "1.835.78"
$ws.Range("E12").Value = "  +0.93%  "

# --- Row 13: WrappedEther ---
$ws.Range("D13").Value = "1.619.42"
$ws.Range("E13").Value = "  +1.29%  "

# --- Row 14: Polkadot ---
$ws.Range("D14").Value = "4.04"
$ws.Range("E14").Value = "  -0.17%  "

# --- Row 15: Polygon ---
$ws.Range("D15").Value = "0.513"
$ws.Range("E15").Value = "  +0.88%  "

# --- Row 16: WrappedBTC ---
$ws.Range("D16").Value = "26.207.89"
$ws.Range("E16").Value = "  +1.63%  "

# --- Row 17: Litecoin ---
$ws.Range("D17").Value = "60.86"
$ws.Range("E17").Value = "  +1.02%  "

# --- Row 18: ShibaInu ---
$ws.Range("E18").Value = "  +2.21%  "

# --- Row 19: Dai ---
$ws.Range("E19").Value = "  -0.46%  "

# --- Row 20: BitcoinCash ---
$ws.Range("D20").Value = "200.25"
$ws.Range("E20").Value = "  +6.19%  "

# --- Row 21: Uniswap ---
$ws.Range("D21").Value = "4.26"
$ws.Range("E21").Value = "  +2.25%  "

# --- Row 22: Avalanche ---
$ws.Range("D22").Value = "9.46"
$ws.Range("E22").Value = "  +1.78%  "

# --- Row 23: Chainlink ---
$ws.Range("D23").Value = "6.03"
$ws.Range("E23").Value = "  +1.72%  "

# --- Row 24: Stellar ---
$ws.Range("D24").Value = "0.133"
$ws.Range("E24").Value = "  +5.49%  "

# --- Row 25: Monero ---
$ws.Range("D25").Value = "142.49"
$ws.Range("E25").Value = "  +0.79%  "

# --- Row 26: Toncoin ---
$ws.Range("E26").Value = "  +2.05%  "

# --- Row 27: BinanceUSD ---
$ws.Range("E27").Value = "  -0.54%  "

# --- Row 28: EthereumClassic ---
$ws.Range("D28").Value = "15.22"
$ws.Range("E28").Value = "  +2.24%  "

# --- Row 29: Cosmos ---
$ws.Range("D29").Value = "6.50"
$ws.Range("E29").Value = "  +0.07%  "

# --- Row 30: PancakeSwap ---
$ws.Range("D30").Value = "1.17"
$ws.Range("E30").Value = "  -0.47%  "

# --- Row 31: Hedera ---
$ws.Range("E31").Value = "  +2.48%  "

# --- Row 32: Filecoin ---
$ws.Range("D32").Value = "3.14"
$ws.Range("E32").Value = "  +2.29%  "

# --- Row 33: InternetComputer(DFINITY) ---
$ws.Range("D33").Value = "3.04"
$ws.Range("E33").Value = "  +2.07%  "

# --- Row 34: LidoDAOToken ---
$ws.Range("D34").Value = "1.51"
$ws.Range("E34").Value = "  +3.44%  "

# --- Row 35: HuobiToken ---
$ws.Range("E35").Value = "  -1.80%  "

# --- Row 36: Maker ---
$ws.Range("D36").Value = "1.110.42"
$ws.Range("E36").Value = "  +1.22%  "

# --- Row 41: ARBITRUM ---
$ws.Range("D41").Value = "0.791"
$ws.Range("E41").Value = "  -0.28%  "

# --- Row 42: TrustWalletToken ---
$ws.Range("D42").Value = "0.794"
$ws.Range("E42").Value = "  +7.48%  "

# --- Row 43: RocketPoolETH ---
$ws.Range("D43").Value = "1.749.08"
$ws.Range("E43").Value = "  +1.02%  "

# --- Row 44: FraxShare ---
$ws.Range("E44").Value = "  +1.65%  "

# --- Row 45: Quant ---
$ws.Range("D45").Value = "93.12"
$ws.Range("E45").Value = "  -2.49%  "

# --- Row 46: RenderToken ---
$ws.Range("E46").Value = "  +9.45%  "

# --- Row 47: BabyDogeCoin ---
$ws.Range("D47").Value = "0.0₆0106"
$ws.Range("E47").Value = "  -5.97%  "

# --- Row 48: Aave ---
$ws.Range("D48").Value = "53.84"
$ws.Range("E48").Value = "  +1.52%  "

# --- Row 49: Cronos ---
$ws.Range("E49").Value = "  +0.01%  "

# --- Row 50: Mantle ---
$ws.Range("E50").Value = "  -0.30%  "

# --- Row 51: USDD ---
$ws.Range("E51").Value = "  -0.34%  "

# --- Rows 37-40: coin rotation (VeChain/PaxDollar/ImmutableX/MXToken) ---
# Row 37 becomes MXToken (was row 40's coin, refreshed numbers)
$ws.Range("B37").Value = "MXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D37").Value = "2.36"
$ws.Range("E37").Value = "  -0.08%  "

# Row 38 becomes VeChain (was row 37's coin, refreshed numbers)
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "0.0153"
$ws.Range("E38").Value = "  +1.55%  "

# Row 39 becomes PaxDollar (was row 38's coin, refreshed numbers)
$ws.Range("B39").Value = "PaxDollar"
$ws.Range("C39").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D39").Value = "1.00"
$ws.Range("E39").Value = "  -0.59%  "

# Row 40 becomes ImmutableX (was row 39's coin, refreshed numbers)
$ws.Range("B40").Value = "ImmutableX"
$ws.Range("C40").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D40").Value = "0.505"
$ws.Range("E40").Value = "  +2.44%  "

